# Applies the "reportes" batch of edits described in the commit:
#   - Clientes: fix a garbled row + append 2 more sample rows
#   - Productos / VentaProductos / ReservasServicios / Facturas:
#       refresh the "Fecha" timestamps and append 3 more sample rows each
#   - Usuarios: append 3 more sample rows

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Clientes
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Clientes")

$ws.Range("A6").Value = 1234
$ws.Range("B6").Value = "migue"
$ws.Range("C6").Value = 123456789

$ws.Range("A11").Value = 1234
$ws.Range("B11").Value = "migue"
$ws.Range("C11").Value = 123456789

$ws.Range("A12").Value = 1234
$ws.Range("B12").Value = "migue"
$ws.Range("C12").Value = 123456789

# ---------------------------------------------------------------------------
# Productos
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Productos")

$ws.Range("H2").Value = "01/06/2024 03:17"
$ws.Range("H3").Value = "01/06/2024 03:19"
$ws.Range("H4").Value = "01/06/2024 03:23"
$ws.Range("H5").Value = "01/06/2024 03:24"
$ws.Range("H6").Value = "01/06/2024 03:25"
$ws.Range("H7").Value = "01/06/2024 03:26"
$ws.Range("H8").Value = "01/06/2024 03:29"
$ws.Range("H9").Value = "01/06/2024 03:33"

$prodRows = @(
    @{ r = 10; h = "01/06/2024 03:36" },
    @{ r = 11; h = "01/06/2024 03:38" },
    @{ r = 12; h = "01/06/2024 03:39" }
)
foreach ($row in $prodRows) {
    $r = $row.r
    $ws.Range("A$r").Value = "REF123"
    $ws.Range("B$r").Value = 1234567890123
    $ws.Range("C$r").Value = "Marca A"
    $ws.Range("D$r").Value = 10
    $ws.Range("E$r").Value = 20
    $ws.Range("F$r").Value = 0
    $ws.Range("G$r").Value = $false
    $ws.Range("H$r").Value = $row.h
}

# ---------------------------------------------------------------------------
# VentaProductos
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("VentaProductos")

$ws.Range("E2").Value = "01/06/2024 03:17"
$ws.Range("E3").Value = "01/06/2024 03:19"
$ws.Range("E4").Value = "01/06/2024 03:23"
$ws.Range("E5").Value = "01/06/2024 03:24"
$ws.Range("E6").Value = "01/06/2024 03:25"
$ws.Range("E7").Value = "01/06/2024 03:26"
$ws.Range("E8").Value = "01/06/2024 03:29"
$ws.Range("E9").Value = "01/06/2024 03:33"

$ventaRows = @(
    @{ r = 10; e = "01/06/2024 03:36" },
    @{ r = 11; e = "01/06/2024 03:38" },
    @{ r = 12; e = "01/06/2024 03:39" }
)
foreach ($row in $ventaRows) {
    $r = $row.r
    $ws.Range("A$r").Value = 123
    $ws.Range("B$r").Value = 1234
    $ws.Range("C$r").Value = "migue"
    $ws.Range("D$r").Value = "Shampoo"
    $ws.Range("E$r").Value = $row.e
    $ws.Range("F$r").Value = 12
    $ws.Range("G$r").Value = 3221
    $ws.Range("H$r").Value = "efectivo"
}

# ---------------------------------------------------------------------------
# ReservasServicios
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ReservasServicios")

$ws.Range("D2").Value = "01/06/2024 03:17"
$ws.Range("D3").Value = "01/06/2024 03:19"
$ws.Range("D4").Value = "01/06/2024 03:23"
$ws.Range("D5").Value = "01/06/2024 03:24"
$ws.Range("D6").Value = "01/06/2024 03:25"
$ws.Range("D7").Value = "01/06/2024 03:26"
$ws.Range("D8").Value = "01/06/2024 03:29"
$ws.Range("D9").Value = "01/06/2024 03:33"

$reservaRows = @(
    @{ r = 10; d = "01/06/2024 03:36" },
    @{ r = 11; d = "01/06/2024 03:38" },
    @{ r = 12; d = "01/06/2024 03:39" }
)
foreach ($row in $reservaRows) {
    $r = $row.r
    $ws.Range("A$r").Value = 12
    $ws.Range("B$r").Value = "mgiue"
    $ws.Range("C$r").Value = 345
    $ws.Range("D$r").Value = $row.d
    $ws.Range("E$r").Value = "28/05/2024 14:30"
}

# ---------------------------------------------------------------------------
# Facturas
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Facturas")

$ws.Range("C2").Value = "01/06/2024 03:17"
$ws.Range("C3").Value = "01/06/2024 03:19"
$ws.Range("C4").Value = "01/06/2024 03:23"
$ws.Range("C5").Value = "01/06/2024 03:24"
$ws.Range("C6").Value = "01/06/2024 03:25"
$ws.Range("C7").Value = "01/06/2024 03:26"
$ws.Range("C8").Value = "01/06/2024 03:29"
$ws.Range("C9").Value = "01/06/2024 03:33"

$facturaRows = @(
    @{ r = 10; c = "01/06/2024 03:36" },
    @{ r = 11; c = "01/06/2024 03:38" },
    @{ r = 12; c = "01/06/2024 03:39" }
)
foreach ($row in $facturaRows) {
    $r = $row.r
    $ws.Range("A$r").Value = 1234
    $ws.Range("B$r").Value = "migue"
    $ws.Range("C$r").Value = $row.c
    $ws.Range("D$r").Value = 3221
    $ws.Range("E$r").Value = 0
    $ws.Range("F$r").Value = 3221
    $ws.Range("G$r").Value = 1
}

# ---------------------------------------------------------------------------
# Usuarios
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Usuarios")

foreach ($r in 11..13) {
    $ws.Range("A$r").Value = 2
    $ws.Range("B$r").Value = "admin"
    $ws.Range("C$r").Value = 12345
    $ws.Range("D$r").Value = 2
}
